$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "69.249.23"
Set-TextValue "E2" "  +4.69%  "
Set-TextValue "D3" "3.542.84"
Set-TextValue "E3" "  +15.43%  "
Set-TextValue "E4" "  -0.16%  "
Set-TextValue "D5" "587.97"
Set-TextValue "E5" "  +2.62%  "
Set-TextValue "D6" "184.42"
Set-TextValue "E6" "  +9.21%  "
Set-TextValue "D7" "3.540.17"
Set-TextValue "E7" "  +15.40%  "
Set-TextValue "D8" "0.998"
Set-TextValue "E8" "  -0.20%  "
Set-TextValue "D9" "0.530"
Set-TextValue "E9" "  +4.15%  "
Set-TextValue "D10" "6.55"
Set-TextValue "E10" "  +4.35%  "
Set-TextValue "D11" "0.158"
Set-TextValue "E11" "  +6.38%  "
Set-TextValue "D12" "0.489"
Set-TextValue "E12" "  +5.02%  "
Set-TextValue "D13" "0.0000250"
Set-TextValue "E13" "  +5.23%  "
Set-TextValue "D14" "38.33"
Set-TextValue "E14" "  +7.77%  "
Set-TextValue "D15" "4.128.61"
Set-TextValue "E15" "  +15.25%  "
Set-TextValue "D16" "69.413.87"
Set-TextValue "E16" "  +4.99%  "
Set-TextValue "B17" "WrappedEther"
Set-TextValue "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "3.542.33"
Set-TextValue "E17" "  +15.22%  "
Set-TextValue "B18" "TRON"
Set-TextValue "C18" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D18" "0.123"
Set-TextValue "E18" "  +1.34%  "
Set-TextValue "D19" "7.43"
Set-TextValue "E19" "  +7.63%  "
Set-TextValue "D20" "16.82"
Set-TextValue "E20" "  +2.63%  "
Set-TextValue "D21" "504.00"
Set-TextValue "E21" "  +4.30%  "
Set-TextValue "D22" "9.06"
Set-TextValue "E22" "  +19.01%  "
Set-TextValue "D23" "0.736"
Set-TextValue "E23" "  +7.88%  "
Set-TextValue "D24" "86.23"
Set-TextValue "E24" "  +4.83%  "
Set-TextValue "D25" "13.34"
Set-TextValue "E25" "  +5.93%  "
Set-TextValue "D26" "2.36"
Set-TextValue "E26" "  +7.77%  "
Set-TextValue "D27" "10.67"
Set-TextValue "E27" "  +5.27%  "
Set-TextValue "E28" "  +0.00%  "
Set-TextValue "D29" "2.51"
Set-TextValue "E29" "  +12.09%  "
Set-TextValue "D30" "8.02"
Set-TextValue "E30" "  +2.32%  "
Set-TextValue "D31" "31.38"
Set-TextValue "E31" "  +13.78%  "
Set-TextValue "B32" "PEPE"
Set-TextValue "C32" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D32" "0.0000109"
Set-TextValue "E32" "  +21.57%  "
Set-TextValue "B33" "PancakeSwap"
Set-TextValue "C33" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D33" "2.71"
Set-TextValue "E33" "  +4.98%  "
Set-TextValue "E34" "  +5.44%  "
Set-TextValue "D35" "0.998"
Set-TextValue "E35" "  -0.21%  "
Set-TextValue "D36" "6.08"
Set-TextValue "E36" "  +9.80%  "
Set-TextValue "E37" "  +7.73%  "
Set-TextValue "D38" "0.330"
Set-TextValue "E38" "  +10.95%  "
Set-TextValue "D39" "2.09"
Set-TextValue "E39" "  +7.39%  "
Set-TextValue "D40" "46.40"
Set-TextValue "E40" "  -1.84%  "
Set-TextValue "D41" "50.58"
Set-TextValue "E41" "  +3.25%  "
Set-TextValue "E42" "  +3.63%  "
Set-TextValue "D43" "8.73"
Set-TextValue "E43" "  +6.59%  "
Set-TextValue "D44" "3.011.84"
Set-TextValue "E44" "  +9.07%  "
Set-TextValue "E45" "  +11.69%  "
Set-TextValue "D46" "398.66"
Set-TextValue "E46" "  +10.54%  "
Set-TextValue "D47" "0.0361"
Set-TextValue "E47" "  +5.91%  "
Set-TextValue "E48" "  +14.12%  "
Set-TextValue "D49" "134.70"
Set-TextValue "E49" "  +0.15%  "
Set-TextValue "E50" "  +0.02%  "
Set-TextValue "D51" "2.44"
Set-TextValue "E51" "  +14.70%  "
